$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.2804546356201172
$ws.Cells.Item(2, 2).Value = 0.4303635954856872
$ws.Cells.Item(2, 3).Value = -0.691750168800354
$ws.Cells.Item(3, 1).Value = 0.1987819671630859
$ws.Cells.Item(3, 2).Value = 0.2879692316055298
$ws.Cells.Item(3, 3).Value = -0.9282988905906676
$ws.Cells.Item(4, 1).Value = 0.2130537033081054
$ws.Cells.Item(4, 2).Value = 0.492556095123291
$ws.Cells.Item(4, 3).Value = -0.9500083923339844
$ws.Cells.Item(5, 1).Value = 5.370540142059326
$ws.Cells.Item(5, 2).Value = -2.71519660949707
$ws.Cells.Item(5, 3).Value = -10.77233219146728
$ws.Cells.Item(6, 1).Value = 1.807073593139648
$ws.Cells.Item(6, 2).Value = 1.692888975143433
$ws.Cells.Item(6, 3).Value = 0.979395866394043
$ws.Cells.Item(7, 1).Value = -1.115023136138916
$ws.Cells.Item(7, 2).Value = -0.1916569471359253
$ws.Cells.Item(7, 3).Value = -1.226519584655761
$ws.Cells.Item(8, 1).Value = -0.6978027820587158
$ws.Cells.Item(8, 2).Value = -0.2956833839416504
$ws.Cells.Item(8, 3).Value = 0.6510324478149414
$ws.Cells.Item(9, 1).Value = -3.345348596572876
$ws.Cells.Item(9, 2).Value = 1.023304224014282
$ws.Cells.Item(9, 3).Value = 2.694005012512207
$ws.Cells.Item(10, 1).Value = -9.996002197265623
$ws.Cells.Item(10, 2).Value = -0.0120857954025268
$ws.Cells.Item(10, 3).Value = 1.461801767349243
$ws.Cells.Item(11, 1).Value = 4.229874610900879
$ws.Cells.Item(11, 2).Value = -2.616225481033325
$ws.Cells.Item(11, 3).Value = -11.63714218139648
$ws.Cells.Item(12, 1).Value = -5.585046291351318
$ws.Cells.Item(12, 2).Value = -1.405134916305542
$ws.Cells.Item(12, 3).Value = 1.352597236633301
$ws.Cells.Item(13, 1).Value = -3.109971046447754
$ws.Cells.Item(13, 2).Value = -0.2804408073425293
$ws.Cells.Item(13, 3).Value = -1.378412246704102
$ws.Cells.Item(14, 1).Value = -6.216216564178467
$ws.Cells.Item(14, 2).Value = 2.080293655395508
$ws.Cells.Item(14, 3).Value = -6.55400562286377
$ws.Cells.Item(15, 1).Value = -2.478143692016602
$ws.Cells.Item(15, 2).Value = -0.5579037666320801
$ws.Cells.Item(15, 3).Value = -6.963788509368896
$ws.Cells.Item(16, 1).Value = 1.977913856506348
$ws.Cells.Item(16, 2).Value = 1.696393609046936
$ws.Cells.Item(16, 3).Value = 4.733741760253906
$ws.Cells.Item(17, 1).Value = -1.786364078521728
$ws.Cells.Item(17, 2).Value = -0.5252545475959778
$ws.Cells.Item(17, 3).Value = -2.082462787628174
$ws.Cells.Item(18, 1).Value = -3.566806316375732
$ws.Cells.Item(18, 2).Value = 0.4948284029960632
$ws.Cells.Item(18, 3).Value = 1.737120628356934
$ws.Cells.Item(19, 1).Value = -6.714832305908203
$ws.Cells.Item(19, 2).Value = 2.814615249633789
$ws.Cells.Item(19, 3).Value = 5.518161296844482
$ws.Cells.Item(20, 1).Value = -4.827847480773926
$ws.Cells.Item(20, 2).Value = -7.360194683074951
$ws.Cells.Item(20, 3).Value = 7.510563373565674
$ws.Cells.Item(21, 1).Value = -1.930130958557129
$ws.Cells.Item(21, 2).Value = 4.015055179595947
$ws.Cells.Item(21, 3).Value = -5.047464370727539
$ws.Cells.Item(22, 1).Value = -4.561958312988281
$ws.Cells.Item(22, 2).Value = 0.1262733936309814
$ws.Cells.Item(22, 3).Value = -1.993290901184082
$ws.Cells.Item(23, 1).Value = -2.729169845581055
$ws.Cells.Item(23, 2).Value = 3.059413433074951
$ws.Cells.Item(23, 3).Value = -4.533473014831543
$ws.Cells.Item(24, 1).Value = -0.4514303207397461
$ws.Cells.Item(24, 2).Value = -0.07753515243530271
$ws.Cells.Item(24, 3).Value = -1.056098580360413
$ws.Cells.Item(25, 1).Value = 1.037992477416992
$ws.Cells.Item(25, 2).Value = -1.273390769958496
$ws.Cells.Item(25, 3).Value = 0.4362349510192871
$ws.Cells.Item(26, 1).Value = 0.0754270553588867
$ws.Cells.Item(26, 2).Value = 1.646718859672546
$ws.Cells.Item(26, 3).Value = 1.695090532302856
$ws.Cells.Item(27, 1).Value = -0.2560558319091797
$ws.Cells.Item(27, 2).Value = 0.3026316165924072
$ws.Cells.Item(27, 3).Value = -0.4233262538909912
$ws.Cells.Item(28, 1).Value = 0.6335611343383789
$ws.Cells.Item(28, 2).Value = 0.8106564879417419
$ws.Cells.Item(28, 3).Value = -1.443797469139099
$ws.Cells.Item(29, 1).Value = 0.09285736083984369
$ws.Cells.Item(29, 2).Value = 0.7357764840126038
$ws.Cells.Item(29, 3).Value = -1.646607518196106
$ws.Cells.Item(30, 1).Value = 0.0882749557495117
$ws.Cells.Item(30, 2).Value = 0.1726978719234466
$ws.Cells.Item(30, 3).Value = -0.9354652166366576
$ws.Cells.Item(31, 1).Value = 0.2656211853027344
$ws.Cells.Item(31, 2).Value = 0.4902379512786865
$ws.Cells.Item(31, 3).Value = -0.8409426212310791

Write-Output "done"
